$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header cells I1 ("I0") and J1 ("IF") ---
# Copy the formatting (bold font, border, centered alignment) from the
# existing header cell H1 so the new headers match the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Fill in the I0 / IF values for each data row (rows 2-24) ---
$iValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 6
    21 = 1
    22 = 5
    23 = 5
    24 = 4
}

$jValues = @{
    2  = 5
    3  = 4
    4  = 5
    5  = 6
    6  = 3
    7  = 4
    8  = 6
    9  = 5
    10 = 6
    11 = 4
    12 = 6
    13 = 6
    14 = 6
    15 = 6
    16 = 3
    17 = 4
    18 = 6
    19 = 5
    20 = 6
    21 = 2
    22 = 6
    23 = 6
    24 = 5
}

for ($row = 2; $row -le 24; $row++) {
    $ws.Cells.Item($row, 9).Value = $iValues[$row]
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
